$d = $word.ActiveDocument
$result = $d.Content.Find.Execute("PhD", $true, $false, $false, $false, $false,
                         $true, 1, $false, "PhD", 2)
Write-Host "Result: $result"
